$wb = $excel.ActiveWorkbook

# zh-cn sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-21 21:06:36"
$wsZhCn.Range("H2").Value = "2016-03-21 21:07:01"

# de-de sheet: update Correspond Handoff Datetime / Correspond Handback DateTime for row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-21 21:06:40"
$wsDeDe.Range("H2").Value = "2016-03-21 21:07:08"
